$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 168-171: new dates + quality/price reshuffle ---

# Row 168 (Calameño)
$ws.Range("D168").Value = 44585
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 600
$ws.Range("K168").Value = 800
$ws.Range("L168").Value = 850
$ws.Range("M168").Value = 825
$ws.Range("P168").Value = 825

# Row 169 (Calameño)
$ws.Range("D169").Value = 44585
$ws.Range("I169").Value = "Segunda"
$ws.Range("J169").Value = 500
$ws.Range("K169").Value = 700
$ws.Range("L169").Value = 750
$ws.Range("M169").Value = 725
$ws.Range("P169").Value = 725

# Row 170 (Tuna)
$ws.Range("D170").Value = 44585
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 600
$ws.Range("K170").Value = 800
$ws.Range("L170").Value = 850
$ws.Range("M170").Value = 825
$ws.Range("P170").Value = 825

# Row 171 (Tuna)
$ws.Range("D171").Value = 44585
$ws.Range("I171").Value = "Segunda"
$ws.Range("J171").Value = 500
$ws.Range("K171").Value = 700
$ws.Range("L171").Value = 750
$ws.Range("M171").Value = 725
$ws.Range("P171").Value = 725

# --- Update existing rows 172-174: shift to new date/grade/region, row174 variety change ---

# Row 172 (Calameño) - was Primera/O'Higgins 44189, now Extra/Maule 44560
$ws.Range("D172").Value = 44560
$ws.Range("I172").Value = "Extra"
$ws.Range("J172").Value = 1500
$ws.Range("K172").Value = 1000
$ws.Range("L172").Value = 1000
$ws.Range("M172").Value = 1000
$ws.Range("O172").Value = "Región del Maule"
$ws.Range("P172").Value = 1000

# Row 173 (Calameño) - was Extra/O'Higgins 44209, now Primera/Maule 44560
$ws.Range("D173").Value = 44560
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 1600
$ws.Range("K173").Value = 700
$ws.Range("L173").Value = 800
$ws.Range("M173").Value = 750
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 750

# Row 174 - was Calameño/Primera/O'Higgins 44209, now Tuna/Extra/Maule 44560
$ws.Range("D174").Value = 44560
$ws.Range("H174").Value = "Tuna"
$ws.Range("I174").Value = "Extra"
$ws.Range("J174").Value = 1500
$ws.Range("K174").Value = 1000
$ws.Range("L174").Value = 1000
$ws.Range("M174").Value = 1000
$ws.Range("O174").Value = "Región del Maule"
$ws.Range("P174").Value = 1000

# --- Row 175: new content (Tuna / Primera / Maule, 44560) replacing the old Calameño/Segunda/O'Higgins row ---
$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 44560
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = 100112027
$ws.Range("G175").Value = "Melón"
$ws.Range("H175").Value = "Tuna"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 1600
$ws.Range("K175").Value = 800
$ws.Range("L175").Value = 900
$ws.Range("M175").Value = 850
$ws.Range("N175").Value = "$/unidad"
$ws.Range("O175").Value = "Región del Maule"
$ws.Range("P175").Value = 850
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = "Hortaliza"

# --- New row 176: Calameño / Primera / O'Higgins, 44189 (old row172 content) ---
$ws.Range("A176").Value = 7
$ws.Range("B176").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C176").Value = "Ñuble"
$ws.Range("D176").Value = 44189
$ws.Range("D176").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E176").Value = 16
$ws.Range("F176").Value = 100112027
$ws.Range("G176").Value = "Melón"
$ws.Range("H176").Value = "Calameño"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 1200
$ws.Range("K176").Value = 800
$ws.Range("L176").Value = 850
$ws.Range("M176").Value = 825
$ws.Range("N176").Value = "$/unidad"
$ws.Range("O176").Value = "Región de O'Higgins"
$ws.Range("P176").Value = 825
$ws.Range("Q176").Value = 1
$ws.Range("R176").Value = "Hortaliza"

# --- New row 177: Calameño / Extra / O'Higgins, 44209 (old row173 content) ---
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44209
$ws.Range("D177").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = 100112027
$ws.Range("G177").Value = "Melón"
$ws.Range("H177").Value = "Calameño"
$ws.Range("I177").Value = "Extra"
$ws.Range("J177").Value = 1200
$ws.Range("K177").Value = 950
$ws.Range("L177").Value = 1000
$ws.Range("M177").Value = 975
$ws.Range("N177").Value = "$/unidad"
$ws.Range("O177").Value = "Región de O'Higgins"
$ws.Range("P177").Value = 975
$ws.Range("Q177").Value = 1
$ws.Range("R177").Value = "Hortaliza"

# --- New row 178: Calameño / Primera / O'Higgins, 44209 (old row174 content) ---
$ws.Range("A178").Value = 7
$ws.Range("B178").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C178").Value = "Ñuble"
$ws.Range("D178").Value = 44209
$ws.Range("D178").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E178").Value = 16
$ws.Range("F178").Value = 100112027
$ws.Range("G178").Value = "Melón"
$ws.Range("H178").Value = "Calameño"
$ws.Range("I178").Value = "Primera"
$ws.Range("J178").Value = 2100
$ws.Range("K178").Value = 750
$ws.Range("L178").Value = 800
$ws.Range("M178").Value = 779
$ws.Range("N178").Value = "$/unidad"
$ws.Range("O178").Value = "Región de O'Higgins"
$ws.Range("P178").Value = 779
$ws.Range("Q178").Value = 1
$ws.Range("R178").Value = "Hortaliza"

# --- New row 179: Calameño / Segunda / O'Higgins, 44209 (old row175 content, renumbered) ---
$ws.Range("A179").Value = 7
$ws.Range("B179").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C179").Value = "Ñuble"
$ws.Range("D179").Value = 44209
$ws.Range("D179").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E179").Value = 16
$ws.Range("F179").Value = 100112027
$ws.Range("G179").Value = "Melón"
$ws.Range("H179").Value = "Calameño"
$ws.Range("I179").Value = "Segunda"
$ws.Range("J179").Value = 1800
$ws.Range("K179").Value = 550
$ws.Range("L179").Value = 600
$ws.Range("M179").Value = 578
$ws.Range("N179").Value = "$/unidad"
$ws.Range("O179").Value = "Región de O'Higgins"
$ws.Range("P179").Value = 578
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = "Hortaliza"
